$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2 (Background Processing) to TRUE
$ws.Range("G2").Value = $true

# Delete row 3 entirely (it duplicated row 2, only differing by the
# Background Processing flag, which is now captured by row 2 itself)
$ws.Rows("3").Delete()

# Select the entire second row, matching the saved selection state
$ws.Range("A2:XFD2").Select()
